$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 141
    3  = 305
    6  = 1847
    10 = 731
    16 = 1205
    17 = 514
    19 = 778
    21 = 408
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
